$wb = $excel.ActiveWorkbook
$wsCore = $wb.Worksheets.Item("Core")
$wsChart = $wb.Worksheets.Item("Test Chart")

# --- Test Chart sheet: rename code value in A2, and match its style (s=2) ---
# Copy the format (style index 2) from Core!A6 before that cell's style gets normalized below.
$wsCore.Range("A6").Copy()
$wsChart.Range("A2").PasteSpecial(-4122)
$wsChart.Range("A2").Value2 = "PatientChartingDate"

# Column A on Test Chart widens from 13.88 to 17.5 (COM ColumnWidth has a constant +5/6 offset
# versus the raw OOXML width attribute, so subtract it to land exactly on 17.5).
$wsChart.Range("A:A").ColumnWidth = 16.666666666666668

# --- Core sheet: the `code` column (A) is rewritten to duplicate the `type` column (B) value ---
$wsCore.Range("A2").Value2 = $wsCore.Range("B2").Value2
$wsCore.Range("A3").Value2 = $wsCore.Range("B3").Value2
$wsCore.Range("A4").Value2 = $wsCore.Range("B4").Value2
$wsCore.Range("A5").Value2 = $wsCore.Range("B5").Value2
$wsCore.Range("A6").Value2 = $wsCore.Range("B6").Value2

# Row 6 (A6, C6, D6) loses its distinct style (s=2) and normalizes to the common style (s=1),
# matching B6's formatting.
$wsCore.Range("B6").Copy()
$wsCore.Range("A6").PasteSpecial(-4122)
$wsCore.Range("C6").PasteSpecial(-4122)
$wsCore.Range("D6").PasteSpecial(-4122)

$excel.CutCopyMode = 0
